$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'42.124.25"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -3.94%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'2.239.35"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -4.44%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  -0.08%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'244.44"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +2.25%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'0.630"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -5.19%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'68.85"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  -5.05%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D9').Value = "'0.553"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -6.74%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'0.0985"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  -1.95%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'59.00"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  -1.34%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'36.01"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  +10.16%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('E13').Value = "'  -2.76%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'6.71"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -7.32%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'2.576.50"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  -4.32%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'14.96"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  -6.70%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'0.861"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  -4.26%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'2.243.61"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -4.01%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'42.081.87"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  -3.83%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'0.0₃0966"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  -6.04%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'6.22"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -6.50%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'73.07"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  -6.62%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'235.50"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -6.17%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('E24').Value = "'  +10.53%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'0.999"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  -0.18%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = "'  -2.62%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'2.47"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  -0.58%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = "'  -1.31%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'9.95"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -4.16%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'172.41"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  -2.45%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'20.47"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  -7.63%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = "'  -3.14%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = "'  -4.97%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'0.0714"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  -3.98%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'5.26"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  -1.21%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'4.70"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -7.15%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'3.78"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  +1.46%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('B38').Value = "'VeChain"
$ws.Range('B38').Style = 'Normal'
$ws.Range('C38').Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range('C38').Style = 'Normal'
$ws.Range('D38').Value = "'0.0287"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  +6.44%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('B39').Value = "'InjectiveProtocol"
$ws.Range('B39').Style = 'Normal'
$ws.Range('C39').Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range('C39').Style = 'Normal'
$ws.Range('D39').Value = "'22.50"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  +20.14%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'2.29"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  -3.07%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'5.88"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  -7.86%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'66.69"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +2.03%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('B43').Value = "'FraxShare"
$ws.Range('B43').Style = 'Normal'
$ws.Range('C43').Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range('C43').Style = 'Normal'
$ws.Range('D43').Value = "'9.24"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +0.36%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('B44').Value = "'FTXToken"
$ws.Range('B44').Style = 'Normal'
$ws.Range('C44').Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range('C44').Style = 'Normal'
$ws.Range('D44').Value = "'5.00"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -14.10%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('E45').Value = "'  -3.09%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('B46').Value = "'Algorand"
$ws.Range('B46').Style = 'Normal'
$ws.Range('C46').Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range('C46').Style = 'Normal'
$ws.Range('D46').Value = "'0.190"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -2.38%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('B47').Value = "'SynthetixNetwork"
$ws.Range('B47').Style = 'Normal'
$ws.Range('C47').Value = "'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range('C47').Style = 'Normal'
$ws.Range('D47').Value = "'4.61"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  +10.48%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = "'  +0.34%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'1.19"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -2.21%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('B50').Value = "'NEARProtocol"
$ws.Range('B50').Style = 'Normal'
$ws.Range('C50').Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range('C50').Style = 'Normal'
$ws.Range('D50').Value = "'2.31"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  -3.49%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('B51').Value = "'HuobiToken"
$ws.Range('B51').Style = 'Normal'
$ws.Range('C51').Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range('C51').Style = 'Normal'
$ws.Range('D51').Value = "'2.81"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  -2.75%  "
$ws.Range('E51').Style = 'Normal'
